$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.689.53'
$ws.Range("E2").Value = '  +1.11%  '

$ws.Range("D3").Value = '1.644.82'
$ws.Range("E3").Value = '  +0.40%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.527'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.63%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.40'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.39%  '

$ws.Range("E9").Value = '  +0.93%  '

$ws.Range("E10").Value = '  +0.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0895'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.53%  '

$ws.Range("D12").Value = '1.877.78'
$ws.Range("E12").Value = '  +0.45%  '

$ws.Range("D13").Value = '1.648.64'
$ws.Range("E13").Value = '  +0.87%  '

$ws.Range("E14").Value = '  +0.64%  '

$ws.Range("E15").Value = '  +0.78%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.75'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.78%  '

$ws.Range("D17").Value = '27.688.29'
$ws.Range("E17").Value = '  +1.23%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '231.27'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.87%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.69'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.72%  '

$ws.Range("E20").Value = '  +0.74%  '

$ws.Range("E21").Value = '  +0.00%  '

$ws.Range("E22").Value = '  -0.36%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.05'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.51%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.85%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.91'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.31%  '

$ws.Range("E26").Value = '  -0.05%  '

$ws.Range("E27").Value = '  -1.18%  '

$ws.Range("E28").Value = '  +0.96%  '

$ws.Range("E29").Value = '  +0.09%  '

$ws.Range("E30").Value = '  +0.82%  '

$ws.Range("E31").Value = '  +0.81%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.32'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.50%  '

$ws.Range("D33").Value = '1.451.89'
$ws.Range("E33").Value = '  +3.12%  '

$ws.Range("E34").Value = '  +0.62%  '

$ws.Range("E35").Value = '  +0.88%  '

$ws.Range("E36").Value = '  -1.01%  '

$ws.Range("E37").Value = '  +1.56%  '

$ws.Range("E38").Value = '  +0.58%  '

$ws.Range("E39").Value = '  +0.52%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.885'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.82%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '71.07'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +10.35%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.04'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.64%  '

$ws.Range("E43").Value = '  +0.05%  '

$ws.Range("E44").Value = '  +2.90%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.47'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.56%  '

$ws.Range("E46").Value = '  +0.70%  '

$ws.Range("D47").Value = '1.787.55'

$ws.Range("E48").Value = '  +5.43%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '85.98'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.44%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0990'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.19%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.78'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.44%  '

